$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New shared string: section header for the "Flink (Stream)" block ---

# --- Row 85: section title (mirrors rows 1 / 26 / 36 / 43 / 55) ---
$ws.Range("A85").Value = "Flink (Stream)"
$ws.Range("A85").NumberFormat = "0.0%"

# --- Row 86: column headers (mirrors rows 2 / 26 / 44 / 56 / etc.) ---
$ws.Range("A86").Value = "CPU0"
$ws.Range("A86").NumberFormat = "0.0%"
$ws.Range("B86").Value = "CPU1"
$ws.Range("B86").NumberFormat = "0.0%"
$ws.Range("C86").Value = "CPU2"
$ws.Range("C86").NumberFormat = "0.0%"
$ws.Range("D86").Value = "CPU3"
$ws.Range("D86").NumberFormat = "0.0%"
$ws.Range("E86").Value = "CPU4"
$ws.Range("E86").NumberFormat = "0.0%"
$ws.Range("F86").Value = "CPU5"
$ws.Range("F86").NumberFormat = "0.0%"
$ws.Range("G86").Value = "CPU6"
$ws.Range("G86").NumberFormat = "0.0%"
$ws.Range("H86").Value = "CPU7"
$ws.Range("H86").NumberFormat = "0.0%"
$ws.Range("I86").Value = "CPU8"
$ws.Range("I86").NumberFormat = "0.0%"
$ws.Range("J86").Value = "CPU9"
$ws.Range("J86").NumberFormat = "0.0%"
$ws.Range("K86").Value = "CPU Avg"
$ws.Range("K86").NumberFormat = "0.0%"
$ws.Range("M86").Value = "Mem0"
$ws.Range("N86").Value = "Mem1"
$ws.Range("O86").Value = "Mem2"
$ws.Range("P86").Value = "Mem3"
$ws.Range("Q86").Value = "Mem4"
$ws.Range("R86").Value = "Mem5"
$ws.Range("S86").Value = "Mem6"
$ws.Range("T86").Value = "Mem7"
$ws.Range("U86").Value = "Mem8"
$ws.Range("V86").Value = "Mem9"
$ws.Range("W86").Value = "Mem Avg"

# --- Rows 87-90: data + per-row formulas (mirrors the "Storm" block, rows 57-84) ---
# Row 87
$ws.Range("A87").Value = 0.003
$ws.Range("A87").NumberFormat = "0.00%"
$ws.Range("B87").Value = 0.002
$ws.Range("B87").NumberFormat = "0.00%"
$ws.Range("C87").Value = 0.002
$ws.Range("C87").NumberFormat = "0.00%"
$ws.Range("D87").Value = 0.002
$ws.Range("D87").NumberFormat = "0.00%"
$ws.Range("E87").Value = 0.002
$ws.Range("E87").NumberFormat = "0.00%"
$ws.Range("F87").Value = 0.002
$ws.Range("F87").NumberFormat = "0.00%"
$ws.Range("G87").Value = 0.003
$ws.Range("G87").NumberFormat = "0.00%"
$ws.Range("H87").Value = 0.002
$ws.Range("H87").NumberFormat = "0.00%"
$ws.Range("I87").Value = 0.002
$ws.Range("I87").NumberFormat = "0.00%"
$ws.Range("J87").Value = 0.002
$ws.Range("J87").NumberFormat = "0.00%"
$ws.Range("K87").Formula = "=AVERAGE(A87:J87)"
$ws.Range("K87").NumberFormat = "0.0%"
$ws.Range("M87").Value = 0.2648
$ws.Range("M87").NumberFormat = "0.00%"
$ws.Range("N87").Value = 0.2691
$ws.Range("N87").NumberFormat = "0.00%"
$ws.Range("O87").Value = 0.2714
$ws.Range("O87").NumberFormat = "0.00%"
$ws.Range("P87").Value = 0.2734
$ws.Range("P87").NumberFormat = "0.00%"
$ws.Range("Q87").Value = 0.2737
$ws.Range("Q87").NumberFormat = "0.00%"
$ws.Range("R87").Value = 0.2759
$ws.Range("R87").NumberFormat = "0.00%"
$ws.Range("S87").Value = 0.278
$ws.Range("S87").NumberFormat = "0.00%"
$ws.Range("T87").Value = 0.2782
$ws.Range("T87").NumberFormat = "0.00%"
$ws.Range("U87").Value = 0.2785
$ws.Range("U87").NumberFormat = "0.00%"
$ws.Range("V87").Value = 0.2818
$ws.Range("V87").NumberFormat = "0.00%"
$ws.Range("W87").Formula = "=AVERAGE(M87:V87)*4000"
$ws.Range("W87").NumberFormat = "0"

# Row 88
$ws.Range("A88").Value = 0.691
$ws.Range("A88").NumberFormat = "0.00%"
$ws.Range("B88").Value = 0.638
$ws.Range("B88").NumberFormat = "0.00%"
$ws.Range("C88").Value = 0.671
$ws.Range("C88").NumberFormat = "0.00%"
$ws.Range("D88").Value = 0.713
$ws.Range("D88").NumberFormat = "0.00%"
$ws.Range("E88").Value = 0.643
$ws.Range("E88").NumberFormat = "0.00%"
$ws.Range("F88").Value = 0.715
$ws.Range("F88").NumberFormat = "0.00%"
$ws.Range("G88").Value = 0.761
$ws.Range("G88").NumberFormat = "0.00%"
$ws.Range("H88").Value = 0.305
$ws.Range("H88").NumberFormat = "0.00%"
$ws.Range("I88").Value = 0.714
$ws.Range("I88").NumberFormat = "0.00%"
$ws.Range("J88").Value = 0.621
$ws.Range("J88").NumberFormat = "0.00%"
$ws.Range("K88").Formula = "=AVERAGE(A88:J88)"
$ws.Range("K88").NumberFormat = "0.0%"
$ws.Range("M88").Value = 0.2648
$ws.Range("M88").NumberFormat = "0.00%"
$ws.Range("N88").Value = 0.2691
$ws.Range("N88").NumberFormat = "0.00%"
$ws.Range("O88").Value = 0.2716
$ws.Range("O88").NumberFormat = "0.00%"
$ws.Range("P88").Value = 0.2734
$ws.Range("P88").NumberFormat = "0.00%"
$ws.Range("Q88").Value = 0.2737
$ws.Range("Q88").NumberFormat = "0.00%"
$ws.Range("R88").Value = 0.2759
$ws.Range("R88").NumberFormat = "0.00%"
$ws.Range("S88").Value = 0.278
$ws.Range("S88").NumberFormat = "0.00%"
$ws.Range("T88").Value = 0.2782
$ws.Range("T88").NumberFormat = "0.00%"
$ws.Range("U88").Value = 0.2787
$ws.Range("U88").NumberFormat = "0.00%"
$ws.Range("V88").Value = 0.282
$ws.Range("V88").NumberFormat = "0.00%"
$ws.Range("W88").Formula = "=AVERAGE(M88:V88)*4000"
$ws.Range("W88").NumberFormat = "0"

# Row 89
$ws.Range("A89").Value = 0.485
$ws.Range("A89").NumberFormat = "0.00%"
$ws.Range("B89").Value = 0.457
$ws.Range("B89").NumberFormat = "0.00%"
$ws.Range("C89").Value = 0.437
$ws.Range("C89").NumberFormat = "0.00%"
$ws.Range("D89").Value = 0.405
$ws.Range("D89").NumberFormat = "0.00%"
$ws.Range("E89").Value = 0.393
$ws.Range("E89").NumberFormat = "0.00%"
$ws.Range("F89").Value = 0.406
$ws.Range("F89").NumberFormat = "0.00%"
$ws.Range("G89").Value = 0.323
$ws.Range("G89").NumberFormat = "0.00%"
$ws.Range("H89").Value = 0.719
$ws.Range("H89").NumberFormat = "0.00%"
$ws.Range("I89").Value = 0.397
$ws.Range("I89").NumberFormat = "0.00%"
$ws.Range("J89").Value = 0.444
$ws.Range("J89").NumberFormat = "0.00%"
$ws.Range("K89").Formula = "=AVERAGE(A89:J89)"
$ws.Range("K89").NumberFormat = "0.0%"
$ws.Range("M89").Value = 0.3035
$ws.Range("M89").NumberFormat = "0.00%"
$ws.Range("N89").Value = 0.3061
$ws.Range("N89").NumberFormat = "0.00%"
$ws.Range("O89").Value = 0.3078
$ws.Range("O89").NumberFormat = "0.00%"
$ws.Range("P89").Value = 0.3101
$ws.Range("P89").NumberFormat = "0.00%"
$ws.Range("Q89").Value = 0.3051
$ws.Range("Q89").NumberFormat = "0.00%"
$ws.Range("R89").Value = 0.3142
$ws.Range("R89").NumberFormat = "0.00%"
$ws.Range("S89").Value = 0.3147
$ws.Range("S89").NumberFormat = "0.00%"
$ws.Range("T89").Value = 0.2957
$ws.Range("T89").NumberFormat = "0.00%"
$ws.Range("U89").Value = 0.3149
$ws.Range("U89").NumberFormat = "0.00%"
$ws.Range("V89").Value = 0.3205
$ws.Range("V89").NumberFormat = "0.00%"
$ws.Range("W89").Formula = "=AVERAGE(M89:V89)*4000"
$ws.Range("W89").NumberFormat = "0"

# Row 90
$ws.Range("A90").Value = 0.003
$ws.Range("A90").NumberFormat = "0.00%"
$ws.Range("B90").Value = 0.003
$ws.Range("B90").NumberFormat = "0.00%"
$ws.Range("C90").Value = 0.002
$ws.Range("C90").NumberFormat = "0.00%"
$ws.Range("D90").Value = 0.003
$ws.Range("D90").NumberFormat = "0.00%"
$ws.Range("E90").Value = 0.002
$ws.Range("E90").NumberFormat = "0.00%"
$ws.Range("F90").Value = 0.003
$ws.Range("F90").NumberFormat = "0.00%"
$ws.Range("G90").Value = 0.002
$ws.Range("G90").NumberFormat = "0.00%"
$ws.Range("H90").Value = 0.062
$ws.Range("H90").NumberFormat = "0.00%"
$ws.Range("I90").Value = 0.002
$ws.Range("I90").NumberFormat = "0.00%"
$ws.Range("J90").Value = 0.002
$ws.Range("J90").NumberFormat = "0.00%"
$ws.Range("K90").Formula = "=AVERAGE(A90:J90)"
$ws.Range("K90").NumberFormat = "0.0%"
$ws.Range("M90").Value = 0.2694
$ws.Range("M90").NumberFormat = "0.00%"
$ws.Range("N90").Value = 0.2716
$ws.Range("N90").NumberFormat = "0.00%"
$ws.Range("O90").Value = 0.2734
$ws.Range("O90").NumberFormat = "0.00%"
$ws.Range("P90").Value = 0.2737
$ws.Range("P90").NumberFormat = "0.00%"
$ws.Range("Q90").Value = 0.2757
$ws.Range("Q90").NumberFormat = "0.00%"
$ws.Range("R90").Value = 0.278
$ws.Range("R90").NumberFormat = "0.00%"
$ws.Range("S90").Value = 0.2782
$ws.Range("S90").NumberFormat = "0.00%"
$ws.Range("T90").Value = 0.317
$ws.Range("T90").NumberFormat = "0.00%"
$ws.Range("U90").Value = 0.282
$ws.Range("U90").NumberFormat = "0.00%"
$ws.Range("V90").Value = 0.2823
$ws.Range("V90").NumberFormat = "0.00%"
$ws.Range("W90").Formula = "=AVERAGE(M90:V90)*4000"
$ws.Range("W90").NumberFormat = "0"

# --- Sheet view: scroll position + active selection, matches the post-edit cursor ---
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("U90").Select()

